$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: Volume/Number and the week-covering date range ---
# "Volume 31   Number  17" -> "...18" (digits 21-22 of the rich-text run)
$ws.Range("A8").Characters(21, 2).Text = "18"
# "Report Covering the Week  4/22/2024  Through  4/28/2024" -> new week
$ws.Range("C9").Characters(27, 9).Text = "4/29/2024"
$ws.Range("C9").Characters(47, 9).Text = "5/5/2024"

# --- Crime Complaints table (rows 14-30): new weekly figures ---

# Row 14: Murder
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = 100
$ws.Range("N14").Value = -33.333333333333

# Row 15: Rape
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = 0
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 80
$ws.Range("L15").Value = 80
$ws.Range("M15").Value = 28.571428571428
$ws.Range("N15").Value = 12.5

# Row 16: Robbery
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 33.870967741935
$ws.Range("L16").Value = 45.614035087719
$ws.Range("M16").Value = 5.063291139240
$ws.Range("N16").Value = -45.394736842105

# Row 17: Fel. Assault
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 6
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 38.095238095238
$ws.Range("I17").Value = 111
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 11
$ws.Range("L17").Value = 46.052631578947
$ws.Range("M17").Value = 79.032258064516
$ws.Range("N17").Value = 48

# Row 18: Burglary
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -7.547169811320
$ws.Range("L18").Value = 32.432432432432
$ws.Range("M18").Value = -42.352941176470
$ws.Range("N18").Value = -83.445945945946

# Row 19: Gr. Larceny
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -43.75
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -37.037037037037
$ws.Range("I19").Value = 207
$ws.Range("J19").Value = 225
$ws.Range("K19").Value = -8
$ws.Range("L19").Value = 9.523809523809
$ws.Range("M19").Value = 55.639097744360
$ws.Range("N19").Value = 48.920863309352

# Row 20: G.L.A.
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -44.444444444444
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 49
$ws.Range("H20").Value = -46.938775510204
$ws.Range("I20").Value = 133
$ws.Range("J20").Value = 179
$ws.Range("K20").Value = -25.698324022346
$ws.Range("L20").Value = 51.136363636363
$ws.Range("M20").Value = 133.333333333333
$ws.Range("N20").Value = -79.969879518072

# Row 21: TOTAL
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -29.729729729729
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = -24.161073825503
$ws.Range("I21").Value = 594
$ws.Range("J21").Value = 625
$ws.Range("K21").Value = -4.96
$ws.Range("L21").Value = 31.415929203539
$ws.Range("M21").Value = 40.425531914893
$ws.Range("N21").Value = -55.572176514584

# Row 22: Transit
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -42.857142857142

# Row 23: Housing
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 60
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = 83.333333333333
$ws.Range("M23").Value = 69.230769230769

# Row 24: Petit Larceny
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -7.142857142857
$ws.Range("F24").Value = 118
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = 9.259259259259
$ws.Range("I24").Value = 486
$ws.Range("J24").Value = 465
$ws.Range("K24").Value = 4.516129032258
$ws.Range("L24").Value = 27.559055118110
$ws.Range("M24").Value = 11.724137931034

# Row 25: Retail Theft
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = -8.571428571428
$ws.Range("I25").Value = 278
$ws.Range("J25").Value = 259
$ws.Range("K25").Value = 7.335907335907
$ws.Range("L25").Value = 47.089947089947

# Row 26: Misd. Assault
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -8.333333333333
$ws.Range("F26").Value = 52
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 23.809523809523
$ws.Range("I26").Value = 193
$ws.Range("J26").Value = 173
$ws.Range("K26").Value = 11.560693641618
$ws.Range("L26").Value = 21.383647798742
$ws.Range("M26").Value = 47.328244274809

# Row 27: UCR Rape*
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = 30
$ws.Range("L27").Value = 8.333333333333

# Row 28: Other Sex Crimes
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 93.333333333333
$ws.Range("L28").Value = 61.111111111111

# Row 29: Shooting Vic.
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -60
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -60

# Row 30: Shooting Inc.
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -50
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = -60
